$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (columns C, D, E) per new contingency results ---

# Row 8
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# Row 13
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- Append new rows 16 and 17 ---

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Copy the style from A15 (bold + border + centered alignment) onto the new
# A16/A17 cells so they match the rest of the "name" column formatting.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
